# The Loan Product short name was missing a hyphen after "295".
# Correct it on both the ProductLoanInput and ProductLoanOutput sheets,
# then leave the ProductLoanOutput sheet active with cell B1 selected
# (mirroring the final view state captured after the edit/verification).

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item(1)   # ProductLoanInput
$wsOutput = $wb.Worksheets.Item(2)   # ProductLoanOutput

$correctedName = "295-MS-EPP-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"

$wsInput.Range("B1").Value = $correctedName
$wsOutput.Range("B1").Value = $correctedName

# Reset selection on the input sheet to B1 (was A6:B6).
$wsInput.Range("B1").Select()

# Make the output sheet the active/selected tab with B1 selected.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
